$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 94
$ws.Range("I9").Value = 73
$ws.Range("J9").Value = 104.5
$ws.Range("K9").Value = 73
$ws.Range("L9").Value = 104.5
$ws.Range("M9").Value = 96
$ws.Range("N9").Value = -442.5
$ws.Range("H15").Value = 1453.9259
$ws.Range("I15").Value = 1453.9259
$ws.Range("K15").Value = 4361.7777
$ws.Range("M15").Value = -4192.7777
$ws.Range("H29").Value = 3671
$ws.Range("I29").Value = 121.333336
$ws.Range("J29").Value = 6713.5713
$ws.Range("K29").Value = 364.000008
$ws.Range("L29").Value = 20140.7139
$ws.Range("M29").Value = -83.00000799999998
$ws.Range("N29").Value = -20702.7139
$ws.Range("H42").Value = 274.77777
$ws.Range("I42").Value = 58
$ws.Range("J42").Value = 2009
$ws.Range("K42").Value = 174
$ws.Range("L42").Value = 6027
$ws.Range("M42").Value = 56
$ws.Range("N42").Value = -6487
$ws.Range("H46").Value = 3017
$ws.Range("I46").Value = 3017
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 9051
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -8932
$ws.Range("N46").ClearContents()
$ws.Range("H60").Value = 3017
$ws.Range("I60").Value = 3017
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 9051
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = -8567
$ws.Range("N60").ClearContents()
$ws.Range("H96").Value = 182374.27
$ws.Range("I96").Value = 182374.27
$ws.Range("K96").Value = 547122.8099999999
$ws.Range("M96").Value = -545749.8099999999
$ws.Range("H100").Value = 3100.3076
$ws.Range("I100").Value = 3073.0908
$ws.Range("K100").Value = 3073.0908
$ws.Range("M100").Value = -2532.0908
$ws.Range("H116").Value = 7691.154
$ws.Range("I116").Value = 4997.75
$ws.Range("K116").Value = 4997.75
$ws.Range("M116").Value = -1555.75
$ws.Range("H138").Value = 3747.6
$ws.Range("J138").Value = 5161.5
$ws.Range("L138").Value = 15484.5
$ws.Range("N138").Value = -25764.5

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1720.7667
$ws.Range("I97").Value = 1162.375
$ws.Range("K97").Value = 1162.375
$ws.Range("M97").Value = -666.375
$ws.Range("H102").Value = 2856.3809
$ws.Range("I102").Value = 2352.9412
$ws.Range("K102").Value = 2352.9412
$ws.Range("M102").Value = -730.9412000000002

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 30000
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("H99").Value = 2354.7778
$ws.Range("I99").Value = 2354.7778
$ws.Range("K99").Value = 2354.7778
$ws.Range("M99").Value = -856.7777999999998

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 5540.3184
$ws.Range("I58").Value = 7126.9375
$ws.Range("K58").Value = 7126.9375
$ws.Range("M58").Value = -6923.9375
$ws.Range("H134").Value = 2898.3635
$ws.Range("I134").Value = 2106.75
$ws.Range("K134").Value = 6320.25
$ws.Range("M134").Value = -3785.25
$ws.Range("H136").Value = 5540.3184
$ws.Range("I136").Value = 7126.9375
$ws.Range("K136").Value = 21380.8125
$ws.Range("M136").Value = -18830.8125

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 285
$ws.Range("I31").Value = 285
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 855
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -567
$ws.Range("N31").ClearContents()
$ws.Range("H137").Value = 4560.9414
$ws.Range("I137").Value = 2762.5
$ws.Range("J137").Value = 4800.7334
$ws.Range("K137").Value = 8287.5
$ws.Range("L137").Value = 14402.2002
$ws.Range("M137").Value = -3187.5
$ws.Range("N137").Value = -24602.2002

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 20000
$ws.Range("J40").Value = 20000
$ws.Range("L40").Value = 20000
$ws.Range("N40").Value = -20302
$ws.Range("H97").Value = 31903.904
$ws.Range("I97").Value = 51021.31
$ws.Range("J97").Value = 838.125
$ws.Range("K97").Value = 51021.31
$ws.Range("L97").Value = 838.125
$ws.Range("M97").Value = -50525.31
$ws.Range("N97").Value = -1830.125
$ws.Range("H130").Value = 70779.5
$ws.Range("J130").Value = 70779.5
$ws.Range("L130").Value = 70779.5
$ws.Range("N130").Value = -80819.5

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1286
$ws.Range("I22").Value = 1533.8
$ws.Range("J22").Value = 1038.2
$ws.Range("K22").Value = 1533.8
$ws.Range("L22").Value = 1038.2
$ws.Range("M22").Value = -1238.8
$ws.Range("N22").Value = -1628.2
$ws.Range("H27").Value = 1286
$ws.Range("I27").Value = 1533.8
$ws.Range("J27").Value = 1038.2
$ws.Range("K27").Value = 1533.8
$ws.Range("L27").Value = 1038.2
$ws.Range("M27").Value = -1426.8
$ws.Range("N27").Value = -1252.2
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").ClearContents()
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()
$ws.Range("H54").Value = 42495
$ws.Range("J54").Value = 42495
$ws.Range("L54").Value = 42495
$ws.Range("N54").Value = -43783
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()
$ws.Range("H58").Value = 11585.857
$ws.Range("I58").Value = 10166.333
$ws.Range("K58").Value = 10166.333
$ws.Range("M58").Value = -9906.333000000001
$ws.Range("H100").Value = 25250.143
$ws.Range("I100").Value = 1551.4615
$ws.Range("J100").Value = 333333
$ws.Range("K100").Value = 1551.4615
$ws.Range("L100").Value = 333333
$ws.Range("M100").Value = -1010.4615
$ws.Range("N100").Value = -334415

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 1000
$ws.Range("I21").Value = 1000
$ws.Range("K21").Value = 1000
$ws.Range("M21").Value = -765
$ws.Range("H35").Value = 1000
$ws.Range("I35").Value = 1000
$ws.Range("K35").Value = 1000
$ws.Range("M35").Value = -710
$ws.Range("H96").Value = 2000
$ws.Range("I96").Value = 2000
$ws.Range("K96").Value = 2000
$ws.Range("M96").Value = -627
$ws.Range("H126").Value = 2784.6428
$ws.Range("I126").Value = 2515.4167
$ws.Range("K126").Value = 7546.250100000001
$ws.Range("M126").Value = -5076.250100000001
$ws.Range("H132").Value = 5122.2964
$ws.Range("I132").Value = 5387.2085
$ws.Range("K132").Value = 16161.6255
$ws.Range("M132").Value = -13631.6255
